{"js": "// POI 4.1.0 -> 5.2.3 bump: three stack-trace line numbers shifted as a\n// consequence of the newer JDK/runtime, inside the single bold/red\n// run of the \"Invalid block\" paragraph.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"ThreadPoolExecutor.java:1130)\", \"ThreadPoolExecutor.java:1136)\"],\n  [\"ThreadPoolExecutor.java:630)\", \"ThreadPoolExecutor.java:635)\"],\n  [\"Thread.java:832)\", \"Thread.java:833)\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of found.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.\n# The stack trace text in the bold/red run shifted three line numbers\n# because of the newer JDK/thread-pool internals.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"ThreadPoolExecutor.java:1130)\", \"ThreadPoolExecutor.java:1136)\"),\n    @(\"ThreadPoolExecutor.java:630)\", \"ThreadPoolExecutor.java:635)\"),\n    @(\"Thread.java:832)\", \"Thread.java:833)\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
